# Add "Area" / "Atotal" columns (G,H) and a small Q/A summary (J,K) to Sheet1,
# mirroring the "add area to Q files stn3" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Headers -----------------------------------------------------------
$ws.Range("G1").Value = "Area"
$ws.Range("H1").Value = "Atotal"
$ws.Range("J1").Value = "Atotal"
$ws.Range("K1").Value = "Qtotal"

# --- Area per segment (G2:G15) -----------------------------------------
# G2 is the special first-segment case (reference depth 0 instead of D1).
$ws.Range("G2").Formula = "=(D2-0)*B2/100"

# G3 is a standalone formula (row before the shared-formula block).
$ws.Range("G3").Formula = "=(D3-D2)*B3/100"

# G4:G15 share one formula pattern, anchored at G4, matching D's own
# shared-formula block (D6:D19) one row ahead of the A/B/C data.
$ws.Range("G4:G15").Formula = "=(D4-D3)*B4/100"

# --- Totals --------------------------------------------------------------
$ws.Range("H2").Formula = "=SUM(G2:G11)"

# --- Summary pair (J2:K2) referencing the totals --------------------------
$ws.Range("J2").Formula = "=H2"
$ws.Range("K2").Formula = "=F2"

# --- View state: match the selection left in the saved workbook -----------
$ws.Range("J2:K2").Select()
